$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "z26"
$ws.Range("B8").Value = "z36"
$ws.Range("C8").Value = "z46"

$ws.Rows.Item(18).Delete()

$null = $ws.Rows.Item(8).Select()
